$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new columns I0 (col I) and IF (col J) with header styled like existing headers
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# row -> [I0, IF] values for rows 2..83
$ijData = @{
    2 = @(7, 7)
    3 = @(9, 9)
    4 = @(6, 6)
    5 = @(7, 7)
    6 = @(8, 8)
    7 = @(7, 7)
    8 = @(7, 7)
    9 = @(8, 8)
    10 = @(7, 7)
    11 = @(7, 7)
    12 = @(6, 6)
    13 = @(7, 7)
    14 = @(8, 8)
    15 = @(7, 7)
    16 = @(7, 7)
    17 = @(7, 7)
    18 = @(7, 7)
    19 = @(7, 7)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(7, 7)
    23 = @(5, 5)
    24 = @(6, 6)
    25 = @(8, 8)
    26 = @(8, 8)
    27 = @(8, 8)
    28 = @(9, 9)
    29 = @(7, 7)
    30 = @(5, 6)
    31 = @(9, 9)
    32 = @(7, 7)
    33 = @(8, 8)
    34 = @(7, 7)
    35 = @(7, 7)
    36 = @(8, 8)
    37 = @(6, 6)
    38 = @(8, 8)
    39 = @(8, 8)
    40 = @(7, 7)
    41 = @(8, 8)
    42 = @(7, 7)
    43 = @(9, 9)
    44 = @(7, 7)
    45 = @(10, 10)
    46 = @(7, 7)
    47 = @(7, 7)
    48 = @(7, 7)
    49 = @(9, 9)
    50 = @(8, 8)
    51 = @(7, 7)
    52 = @(10, 10)
    53 = @(8, 8)
    54 = @(4, 4)
    55 = @(7, 7)
    56 = @(8, 8)
    57 = @(5, 5)
    58 = @(8, 8)
    59 = @(6, 6)
    60 = @(7, 7)
    61 = @(11, 11)
    62 = @(4, 4)
    63 = @(8, 8)
    64 = @(7, 7)
    65 = @(7, 7)
    66 = @(6, 6)
    67 = @(6, 6)
    68 = @(8, 8)
    69 = @(8, 8)
    70 = @(6, 7)
    71 = @(7, 7)
    72 = @(8, 8)
    73 = @(10, 10)
    74 = @(9, 9)
    75 = @(9, 9)
    76 = @(8, 8)
    77 = @(7, 7)
    78 = @(4, 5)
    79 = @(5, 5)
    80 = @(7, 7)
    81 = @(5, 5)
    82 = @(5, 5)
    83 = @(4, 4)
}

foreach ($r in $ijData.Keys) {
    $vals = $ijData[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
